$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (Jordan Poole @ row 19 no longer needed there; moves to row 13)
$ws.Rows.Item(19).Delete()

# Update rows 2-18 with the reshuffled roster data
$ws.Range("A2").Value = "Keyonte George"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Utah Jazz"
$ws.Range("A3").Value = "Dalano Banton"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Portland Trail Blazers"
$ws.Range("A4").Value = "Tobias Harris"
$ws.Range("B4").Value = "SF,PF"
$ws.Range("C4").Value = "Detroit Pistons"
$ws.Range("A5").Value = "Lauri Markkanen"
$ws.Range("B5").Value = "SF,PF"
$ws.Range("C5").Value = "Utah Jazz"
$ws.Range("A6").Value = "Christian Braun"
$ws.Range("B6").Value = "SG,SF"
$ws.Range("C6").Value = "Denver Nuggets"
$ws.Range("A7").Value = "Zach LaVine"
$ws.Range("B7").Value = "SG,SF"
$ws.Range("C7").Value = "Chicago Bulls"
$ws.Range("A8").Value = "John Collins"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Utah Jazz"
$ws.Range("A9").Value = "Jalen Williams"
$ws.Range("B9").Value = "SG,SF,PF,C"
$ws.Range("C9").Value = "Oklahoma City Thunder"
$ws.Range("A10").Value = "Shai Gilgeous-Alexander"
$ws.Range("B10").Value = "PG"
$ws.Range("C10").Value = "Oklahoma City Thunder"
$ws.Range("A11").Value = "Dennis Schröder"
$ws.Range("B11").Value = "PG"
$ws.Range("C11").Value = "Brooklyn Nets"
$ws.Range("A12").Value = "Kyrie Irving"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "Dallas Mavericks"
$ws.Range("A13").Value = "Jordan Poole"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Washington Wizards"
$ws.Range("A14").Value = "Brandon Boston Jr."
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "New Orleans Pelicans"
$ws.Range("A15").Value = "RJ Barrett"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "Toronto Raptors"
$ws.Range("A16").Value = "Jimmy Butler"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Miami Heat"
$ws.Range("A17").Value = "CJ McCollum"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "New Orleans Pelicans"
$ws.Range("A18").Value = "Joel Embiid"
$ws.Range("B18").Value = "C"
$ws.Range("C18").Value = "Philadelphia 76ers"
